# Updates the cryptos list (prices and 1h volume %) to match the
# latest scrape, and re-orders three coins whose rank changed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.858.63"
$ws.Range("E2").Value = "  -5.86%  "
$ws.Range("D3").Value = "2.232.91"
$ws.Range("E3").Value = "  -6.68%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.76"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -12.48%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.571"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -9.88%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.561"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -10.26%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.36"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -11.21%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0837"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -9.65%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.60"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -12.56%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.877"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -12.62%  "
$ws.Range("D15").Value = "2.571.79"
$ws.Range("E15").Value = "  -6.65%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.84"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -12.75%  "
$ws.Range("D17").Value = "2.243.49"
$ws.Range("E17").Value = "  -6.81%  "
$ws.Range("D18").Value = "42.646.66"
$ws.Range("E18").Value = "  -6.22%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.34"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.95%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.63"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -11.29%  "
$ws.Range("D21").Value = "0.0₃0950"
$ws.Range("E21").Value = "  -12.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.27"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -7.72%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "64.99"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -12.92%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "234.74"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -11.19%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.13"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -9.36%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.14%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.25"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -9.43%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.62"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -14.35%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.17"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -7.78%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0890"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -8.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.59"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -9.35%  "
$ws.Range("B32").Value = "Monero"
$ws.Range("C32").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "162.27"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.89%  "
$ws.Range("B33").Value = "InjectiveProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "33.58"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -14.73%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.76"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.31%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.06"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.38%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.121"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -8.62%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.88"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.27%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.40"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -10.82%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.106"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -10.46%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.64"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -11.60%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0322"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -10.96%  "
$ws.Range("E42").Value = "  +0.16%  "
$ws.Range("D43").Value = "1.817.11"
$ws.Range("E43").Value = "  +8.92%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "90.16"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -11.13%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.13"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -10.69%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.207"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -14.36%  "
$ws.Range("B47").Value = "THORChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.37"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.10%  "
$ws.Range("B48").Value = "MultiversX"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "61.11"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -14.70%  "
$ws.Range("B49").Value = "ordi"
$ws.Range("C49").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "76.55"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -12.60%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.69"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -8.41%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "102.35"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -11.44%  "
